$wb = $excel.ActiveWorkbook

# Region_Code value changes from ND05 to ND01 on the WMT_Extract sheet (rows 2 & 3)
$wsExtract = $wb.Worksheets.Item("WMT_Extract")
$wsExtract.Range("C2:C3").Value = "ND01"

# Inst_Reports: fill in the previously-blank G2 cell to match F2 (no blank columns)
$wsInst = $wb.Worksheets.Item("Inst_Reports")
$wsInst.Range("G2").Value = 2

# Move the selection/active-cell on WMT_Extract
[void]$wsExtract.Range("AE2").Select()

# Make Inst_Reports the active sheet with its own selection
[void]$wsInst.Select()
[void]$wsInst.Range("E14").Select()
